$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-number-looking string as TEXT (matches original inlineStr cells),
# using a leading apostrophe to force text entry, then resetting the style so no
# extra quote-prefix / number-format style gets attached to the cell.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '62.795.31'
$ws.Range("E2").Value = '  +2.25%  '
Set-TextValue $ws.Range("D3") '3.030.12'
$ws.Range("E3").Value = '  +1.41%  '
Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue $ws.Range("D5") '595.33'
$ws.Range("E5").Value = '  +1.24%  '
Set-TextValue $ws.Range("D6") '151.94'
$ws.Range("E6").Value = '  +6.01%  '
$ws.Range("E7").Value = '  -0.04%  '
Set-TextValue $ws.Range("D8") '3.024.34'
$ws.Range("E8").Value = '  +1.26%  '
Set-TextValue $ws.Range("D9") '0.518'
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("E10").Value = '  +8.61%  '
Set-TextValue $ws.Range("D11") '0.151'
$ws.Range("E11").Value = '  +4.12%  '
Set-TextValue $ws.Range("D12") '0.461'
$ws.Range("E12").Value = '  +0.01%  '
Set-TextValue $ws.Range("D13") '0.0000233'
$ws.Range("E13").Value = '  +2.71%  '
Set-TextValue $ws.Range("D14") '34.97'
$ws.Range("E14").Value = '  +2.16%  '
$ws.Range("E15").Value = '  +2.21%  '
Set-TextValue $ws.Range("D16") '3.530.94'
$ws.Range("E16").Value = '  +1.36%  '
Set-TextValue $ws.Range("D17") '62.767.23'
$ws.Range("E17").Value = '  +2.22%  '
Set-TextValue $ws.Range("D18") '7.02'
$ws.Range("E18").Value = '  -0.34%  '
Set-TextValue $ws.Range("D19") '3.026.90'
$ws.Range("E19").Value = '  +1.33%  '
Set-TextValue $ws.Range("D20") '452.29'
$ws.Range("E20").Value = '  +0.17%  '
Set-TextValue $ws.Range("D21") '14.20'
$ws.Range("E21").Value = '  +1.02%  '
Set-TextValue $ws.Range("D22") '0.692'
$ws.Range("E22").Value = '  +0.98%  '
Set-TextValue $ws.Range("D23") '7.47'
$ws.Range("E23").Value = '  +1.61%  '
Set-TextValue $ws.Range("D24") '83.27'
$ws.Range("E24").Value = '  +1.81%  '
Set-TextValue $ws.Range("D25") '2.27'
$ws.Range("E25").Value = '  +4.69%  '
Set-TextValue $ws.Range("D26") '10.94'
$ws.Range("E26").Value = '  +7.09%  '
Set-TextValue $ws.Range("D27") '12.17'
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("E28").Value = '  +0.00%  '
Set-TextValue $ws.Range("D29") '7.52'
$ws.Range("E29").Value = '  +8.39%  '
Set-TextValue $ws.Range("D30") '2.71'
$ws.Range("E30").Value = '  +2.01%  '
Set-TextValue $ws.Range("D31") '1.00'
$ws.Range("E31").Value = '  -0.04%  '
Set-TextValue $ws.Range("D32") '2.20'
$ws.Range("E32").Value = '  +7.44%  '
Set-TextValue $ws.Range("D33") '27.53'
$ws.Range("E33").Value = '  +1.38%  '
Set-TextValue $ws.Range("D34") '0.110'
$ws.Range("E34").Value = '  +3.22%  '
Set-TextValue $ws.Range("D35") '0.0₃0847'
$ws.Range("E35").Value = '  +5.76%  '
$ws.Range("E36").Value = '  +2.03%  '
Set-TextValue $ws.Range("D37") '5.89'
$ws.Range("E37").Value = '  +2.64%  '
Set-TextValue $ws.Range("D38") '3.09'
$ws.Range("E38").Value = '  +9.80%  '
Set-TextValue $ws.Range("D41") '9.08'
$ws.Range("E41").Value = '  -0.93%  '
Set-TextValue $ws.Range("D42") '0.125'
$ws.Range("E42").Value = '  +3.17%  '
Set-TextValue $ws.Range("D43") '0.298'
$ws.Range("E43").Value = '  +12.98%  '
Set-TextValue $ws.Range("D44") '42.21'
$ws.Range("E44").Value = '  +10.16%  '
Set-TextValue $ws.Range("D45") '390.01'
$ws.Range("E45").Value = '  -3.08%  '
Set-TextValue $ws.Range("D46") '0.0356'
$ws.Range("E46").Value = '  +1.09%  '
Set-TextValue $ws.Range("D47") '2.736.99'
$ws.Range("E47").Value = '  +0.96%  '
Set-TextValue $ws.Range("D48") '131.63'
$ws.Range("E48").Value = '  -0.81%  '
Set-TextValue $ws.Range("D50") '2.20'
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("E51").Value = '  +1.36%  '

# Rows 39/40 swap places: Stacks <-> OKB
Set-TextValue $ws.Range("B39") 'OKB'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D39") '50.64'
$ws.Range("E39").Value = '  +1.09%  '
Set-TextValue $ws.Range("B40") 'Stacks'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D40") '2.07'
$ws.Range("E40").Value = '  -0.19%  '
